$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 125: Anomaly no noise / i6vtqbf9 / Training phase / round 1
$ws.Range("A125").Value = "Anomaly no noise"
$ws.Range("B125").Value = "i6vtqbf9"
$ws.Range("C125").Value = "Training phase"
$ws.Range("D125").Value = 1
$ws.Range("E125").Value = "['Purple', 'Green']"
$ws.Range("F125").Value = "[['Red', 'Blue'], ['Blue', 'Yellow']]"
$ws.Range("G125").Value = "[None, None]"
$ws.Range("H125").Value = "['8', '8']"

# Row 126: No Anomaly no noise / mni5vgdn / Training phase / round 1
$ws.Range("A126").Value = "No Anomaly no noise"
$ws.Range("B126").Value = "mni5vgdn"
$ws.Range("C126").Value = "Training phase"
$ws.Range("D126").Value = 1
$ws.Range("E126").Value = "['Purple', 'Green']"
$ws.Range("F126").Value = "[['Blue', 'Red'], ['Blue', 'Yellow']]"
$ws.Range("G126").Value = "[None, None]"
$ws.Range("H126").Value = "['8', '8']"

# Row 127: No Anomaly no noise / mni5vgdn / Training phase / round 2
$ws.Range("A127").Value = "No Anomaly no noise"
$ws.Range("B127").Value = "mni5vgdn"
$ws.Range("C127").Value = "Training phase"
$ws.Range("D127").Value = 2
$ws.Range("E127").Value = "['Green', 'Purple']"
$ws.Range("F127").Value = "[['Yellow', 'Blue'], ['Blue', 'Red']]"
$ws.Range("G127").Value = "[None, None]"
$ws.Range("H127").Value = "['8', '8']"

# Row 128: Anomaly no noise / i6vtqbf9 / Training phase / round 2
$ws.Range("A128").Value = "Anomaly no noise"
$ws.Range("B128").Value = "i6vtqbf9"
$ws.Range("C128").Value = "Training phase"
$ws.Range("D128").Value = 2
$ws.Range("E128").Value = "['Green', 'Purple']"
$ws.Range("F128").Value = "[['Blue', 'Yellow'], ['Red', 'Blue']]"
$ws.Range("G128").Value = "[None, None]"
$ws.Range("H128").Value = "['8', '8']"

# Row 129: Anomaly no noise / i6vtqbf9 / Training phase / round 3
$ws.Range("A129").Value = "Anomaly no noise"
$ws.Range("B129").Value = "i6vtqbf9"
$ws.Range("C129").Value = "Training phase"
$ws.Range("D129").Value = 3
$ws.Range("E129").Value = "['Purple', 'Green', 'Green']"
$ws.Range("F129").Value = "[['Red', ''], ['Blue', ''], ['Blue', '']]"
$ws.Range("G129").Value = "[None, None, None]"
$ws.Range("H129").Value = "['2', '5', '5']"

# Row 130: No Anomaly no noise / mni5vgdn / Training phase / round 3
$ws.Range("A130").Value = "No Anomaly no noise"
$ws.Range("B130").Value = "mni5vgdn"
$ws.Range("C130").Value = "Training phase"
$ws.Range("D130").Value = 3
$ws.Range("E130").Value = "['Purple', 'Green', 'Green']"
$ws.Range("F130").Value = "[['Red', ''], ['Blue', ''], ['Blue', '']]"
$ws.Range("G130").Value = "[None, None, None]"
$ws.Range("H130").Value = "['2', '5', '5']"

# Row 131: Anomaly no noise / i6vtqbf9 / Training phase / round 4
$ws.Range("A131").Value = "Anomaly no noise"
$ws.Range("B131").Value = "i6vtqbf9"
$ws.Range("C131").Value = "Training phase"
$ws.Range("D131").Value = 4
$ws.Range("E131").Value = "['Green', 'Purple', 'Purple']"
$ws.Range("F131").Value = "[['Red', ''], ['Red', ''], ['Red', '']]"
$ws.Range("G131").Value = "[None, None, None]"
$ws.Range("H131").Value = "['2', '2', '2']"

# Row 132: No Anomaly no noise / mni5vgdn / Training phase / round 4
$ws.Range("A132").Value = "No Anomaly no noise"
$ws.Range("B132").Value = "mni5vgdn"
$ws.Range("C132").Value = "Training phase"
$ws.Range("D132").Value = 4
$ws.Range("E132").Value = "['Green', 'Purple', 'Purple']"
$ws.Range("F132").Value = "[['Red', ''], ['Red', ''], ['Red', '']]"
$ws.Range("G132").Value = "[None, None, None]"
$ws.Range("H132").Value = "['2', '2', '2']"

# Row 133: Anomaly no noise / i6vtqbf9 / Test 1 / round 1
$ws.Range("A133").Value = "Anomaly no noise"
$ws.Range("B133").Value = "i6vtqbf9"
$ws.Range("C133").Value = "Test 1"
$ws.Range("D133").Value = 1
$ws.Range("E133").Value = "['Green', 'Yellow', 'Purple', 'Red', 'Orange', 'Blue']"
$ws.Range("F133").Value = "[['Red', ''], ['Red', ''], ['Red', ''], ['Red', ''], ['Red', ''], ['Red', '']]"
$ws.Range("G133").Value = "[None, None, None, None, None, None]"
$ws.Range("H133").Value = "['0', '0', '0', '0', '0', '0']"

# Row 134: No Anomaly no noise / mni5vgdn / Test 1 / round 1
$ws.Range("A134").Value = "No Anomaly no noise"
$ws.Range("B134").Value = "mni5vgdn"
$ws.Range("C134").Value = "Test 1"
$ws.Range("D134").Value = 1
$ws.Range("E134").Value = "['Green', 'Yellow', 'Purple', 'Red', 'Orange', 'Blue']"
$ws.Range("F134").Value = "[['Red', ''], ['Red', ''], ['Red', ''], ['Red', ''], ['Red', ''], ['Red', '']]"
$ws.Range("G134").Value = "[None, None, None, None, None, None]"
$ws.Range("H134").Value = "['0', '0', '0', '0', '0', '0']"
$ws.Range("I134").NumberFormat = "@"
$ws.Range("I134").Value = "0.50"

# Row 135: Anomaly no noise / i6vtqbf9 / Exploration phase / round 1
$ws.Range("A135").Value = "Anomaly no noise"
$ws.Range("B135").Value = "i6vtqbf9"
$ws.Range("C135").Value = "Exploration phase"
$ws.Range("D135").Value = 1
$ws.Range("E135").Value = "['Orange', 'Purple']"
$ws.Range("F135").Value = "[['Red', ''], ['Red', '']]"
$ws.Range("G135").Value = "[None, None]"
$ws.Range("H135").Value = "['2', '2']"

# Row 136: Anomaly no noise / i6vtqbf9 / Exploration phase / round 2
$ws.Range("A136").Value = "Anomaly no noise"
$ws.Range("B136").Value = "i6vtqbf9"
$ws.Range("C136").Value = "Exploration phase"
$ws.Range("D136").Value = 2
$ws.Range("E136").Value = "['Orange', 'Green']"
$ws.Range("F136").Value = "[['Red', ''], ['Red', '']]"
$ws.Range("G136").Value = "[None, None]"
$ws.Range("H136").Value = "['2', '2']"

# Row 137: Anomaly no noise / i6vtqbf9 / Exploration phase / round 3
$ws.Range("A137").Value = "Anomaly no noise"
$ws.Range("B137").Value = "i6vtqbf9"
$ws.Range("C137").Value = "Exploration phase"
$ws.Range("D137").Value = 3
$ws.Range("E137").Value = "['Green', 'Red', 'Green']"
$ws.Range("F137").Value = "[['Red', ''], ['Blue', ''], ['Blue', '']]"
$ws.Range("G137").Value = "[None, None, None]"
$ws.Range("H137").Value = "['2', '5', '5']"

# Row 138: Anomaly no noise / i6vtqbf9 / Exploration phase / round 4
$ws.Range("A138").Value = "Anomaly no noise"
$ws.Range("B138").Value = "i6vtqbf9"
$ws.Range("C138").Value = "Exploration phase"
$ws.Range("D138").Value = 4
$ws.Range("E138").Value = "['Blue', 'Purple', 'Purple']"
$ws.Range("F138").Value = "[['Red', ''], ['Blue', ''], ['Red', '']]"
$ws.Range("G138").Value = "[None, None, None]"
$ws.Range("H138").Value = "['2', '5', '2']"

# Row 139: Anomaly no noise / i6vtqbf9 / Exploration phase / round 5
$ws.Range("A139").Value = "Anomaly no noise"
$ws.Range("B139").Value = "i6vtqbf9"
$ws.Range("C139").Value = "Exploration phase"
$ws.Range("D139").Value = 5
$ws.Range("E139").Value = "['Purple', 'Green', 'Yellow']"
$ws.Range("F139").Value = "[['Blue', ''], ['Red', ''], ['Red', '']]"
$ws.Range("G139").Value = "[None, None, None]"
$ws.Range("H139").Value = "['5', '2', '2']"
$ws.Range("I139").NumberFormat = "@"
$ws.Range("I139").Value = "0.88"

# Row 140: No Anomaly noisy / anucg6jd / Training phase / round 1
$ws.Range("A140").Value = "No Anomaly noisy"
$ws.Range("B140").Value = "anucg6jd"
$ws.Range("C140").Value = "Training phase"
$ws.Range("D140").Value = 1
$ws.Range("E140").Value = "['Purple', 'Green']"
$ws.Range("F140").Value = "[['Red', 'Red'], ['Blue', 'Blue']]"
$ws.Range("G140").Value = "[{'index': 0, 'type': 'increase', 'amount': 0.2, 'before': 0.5, 'after': 0.7}, None]"
$ws.Range("H140").Value = "['7', '10']"
$ws.Range("I140").NumberFormat = "@"
$ws.Range("I140").Value = "0.17"

# Row 141: Anomaly noisy / tbcsissw / Training phase / round 1
$ws.Range("A141").Value = "Anomaly noisy"
$ws.Range("B141").Value = "tbcsissw"
$ws.Range("C141").Value = "Training phase"
$ws.Range("D141").Value = 1
$ws.Range("E141").Value = "['Purple', 'Green']"
$ws.Range("F141").Value = "[['Red', 'Red'], ['Blue', 'Blue']]"
$ws.Range("G141").Value = "[None, {'index': 1, 'type': 'decrease', 'amount': -0.2, 'before': 1.0, 'after': 0.8}]"
$ws.Range("H141").Value = "['5', '8']"
$ws.Range("I141").NumberFormat = "@"
$ws.Range("I141").Value = "0.13"
